$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Typography")
$ws2 = $wb.Worksheets.Item("Translation")

# Sheet1 (Typography): set I4 ("Wildcard Ranges" column for first font row) to "0-9"
$ws1.Range("I4").Value = "0-9"

# Sheet2 (Translation): update row 5
$ws2.Range("D5").Value = "Right"
$ws2.Range("F5").Value = "<value>"

# Sheet2 (Translation): add new row 6
$ws2.Range("B6").Value = "SingleUseId4"
$ws2.Range("C6").Value = "Default"
$ws2.Range("D6").Value = "Right"
$ws2.Range("E6").Value = "LTR"
$ws2.Range("F6").Value = "<value>"
